$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora): all rows 2-51 change from 21 to 22
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "22"
$ws.Range("G2:G51").Style = "Normal"

# Column D (Price) and E (Volume(1h)) per-row updates
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.89%"
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.52%"
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.084"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.00%"
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.12%"
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.353"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.89%"
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.892"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.70%"
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.193"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.10%"
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.014"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-9.07%"
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.53%"
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-13.24%"
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1863"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.84%"
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08686"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.29%"
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03391"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.53%"
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09697"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.94%"
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001364"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.54%"
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005964"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.55%"
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.617"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.82%"
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.06%"
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1274"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.67%"
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.016"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.89%"
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2596"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.01%"
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,175.32%"
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04320"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.78%"
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.22%"
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004211"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-12.30%"
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.88%"
$ws.Range("E27").Style = "Normal"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02198"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.13%"
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04904"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.56%"
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007548"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.63%"
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.18%"
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.18%"
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002063"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.96%"
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008498"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.02%"
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006554"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.91%"
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.25%"
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.98%"
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-22.88%"
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.25%"
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.25%"
$ws.Range("E51").Style = "Normal"
